# Update market-price derived columns (currentAveragePrice* / LevePrice* / LeveProfit*)
# on the per-job Leve profit sheets, per latest market data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 2690
$ws.Range("J18").Value = 5300
$ws.Range("L18").Value = 5300
$ws.Range("N18").Value = -5868
# Row 29: Dripping with Venom
$ws.Range("H29").Value = 6249.5
$ws.Range("I29").Value = 3000
$ws.Range("K29").Value = 9000
$ws.Range("M29").Value = -8719
# Row 38: Just Give Him a Serum
$ws.Range("H38").Value = 700.3
$ws.Range("I38").Value = 700.3
$ws.Range("K38").Value = 2100.9
$ws.Range("M38").Value = -1728.9
# Row 39: Riches' Brew
$ws.Range("H39").Value = 189.89655
$ws.Range("I39").Value = 82.35714
$ws.Range("J39").Value = 290.26666
$ws.Range("K39").Value = 247.07142
$ws.Range("L39").Value = 870.79998
$ws.Range("M39").Value = 48.92858000000001
$ws.Range("N39").Value = -1462.79998
# Row 42: Eye of the Beholder
$ws.Range("H42").Value = 486
$ws.Range("I42").Value = 324
$ws.Range("J42").Value = 648
$ws.Range("K42").Value = 972
$ws.Range("L42").Value = 1944
$ws.Range("M42").Value = -742
$ws.Range("N42").Value = -2404
# Row 109: A Time for Peace
$ws.Range("H109").Value = 89330
$ws.Range("J109").Value = 89330
$ws.Range("L109").Value = 89330
$ws.Range("N109").Value = -92104
# Row 112: Making Ends Meet
$ws.Range("H112").Value = 2242.4167
$ws.Range("J112").Value = 2264.5454
$ws.Range("L112").Value = 6793.6362
$ws.Range("N112").Value = -9009.636200000001
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1273.68
$ws.Range("I132").Value = 724.0222
$ws.Range("K132").Value = 2172.0666
$ws.Range("M132").Value = 357.9333999999999
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2721.6938
$ws.Range("J138").Value = 2972.6099
$ws.Range("L138").Value = 8917.8297
$ws.Range("N138").Value = -19197.8297

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 18520140
$ws.Range("I45").Value = 20834828
$ws.Range("K45").Value = 20834828
$ws.Range("M45").Value = -20834451
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2863.8667
$ws.Range("I122").Value = 2101.238
$ws.Range("K122").Value = 6303.714
$ws.Range("M122").Value = -3853.714
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3358.3235
$ws.Range("I132").Value = 2442.5
$ws.Range("J132").Value = 10227
$ws.Range("K132").Value = 7327.5
$ws.Range("L132").Value = 30681
$ws.Range("M132").Value = -4797.5
$ws.Range("N132").Value = -35741

$ws = $wb.Worksheets.Item("BSM")
# Row 5: Axe Me Anything
$ws.Range("H5").Value = 1999.5
$ws.Range("J5").Value = 2333
$ws.Range("L5").Value = 2333
$ws.Range("N5").Value = -2559
# Row 94: High Steal
$ws.Range("H94").Value = 1659.9131
$ws.Range("I94").Value = 607.0769
$ws.Range("K94").Value = 607.0769
$ws.Range("M94").Value = -156.0769
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2043.0435
$ws.Range("I105").Value = 1392.5
$ws.Range("J105").Value = 2543.4614
$ws.Range("K105").Value = 1392.5
$ws.Range("L105").Value = 2543.4614
$ws.Range("M105").Value = 354.5
$ws.Range("N105").Value = -6037.4614
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 236666.16
$ws.Range("I134").Value = 1500.5143
$ws.Range("K134").Value = 4501.5429
$ws.Range("M134").Value = -1966.5429

$ws = $wb.Worksheets.Item("CRP")
# Row 28: Militia on My Mind
$ws.Range("H28").Value = 11412
$ws.Range("J28").Value = 11412
$ws.Range("L28").Value = 11412
$ws.Range("N28").Value = -11902
# Row 60: Bowing to Greater Power
$ws.Range("H60").Value = 8249.75
$ws.Range("I60").Value = 8249.75
$ws.Range("K60").Value = 8249.75
$ws.Range("M60").Value = -7738.75
# Row 75: The Darkest Hearth
$ws.Range("H75").Value = 112989
$ws.Range("J75").Value = 112989
$ws.Range("L75").Value = 112989
$ws.Range("N75").Value = -114985
# Row 78: Fruit of the Loom (L)
$ws.Range("H78").Value = 112989
$ws.Range("J78").Value = 112989
$ws.Range("L78").Value = 338967
$ws.Range("N78").Value = -348951
# Row 86: Birch, Please
$ws.Range("H86").Value = 87817.414
$ws.Range("I86").Value = 5533.6665
$ws.Range("J86").Value = 170101.17
$ws.Range("K86").Value = 5533.6665
$ws.Range("L86").Value = 170101.17
$ws.Range("M86").Value = -4410.6665
$ws.Range("N86").Value = -172347.17
# Row 89: Built This City on Blocks and Soul (L)
$ws.Range("H89").Value = 87817.414
$ws.Range("I89").Value = 5533.6665
$ws.Range("J89").Value = 170101.17
$ws.Range("K89").Value = 27668.3325
$ws.Range("L89").Value = 850505.8500000001
$ws.Range("M89").Value = -22052.3325
$ws.Range("N89").Value = -861737.8500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service
$ws.Range("H107").Value = 1771.125
$ws.Range("J107").Value = 1771.125
$ws.Range("L107").Value = 5313.375
$ws.Range("N107").Value = -9153.375
# Row 122: Salt of the North
$ws.Range("H122").Value = 2639.3462
$ws.Range("I122").Value = 700.3333
$ws.Range("K122").Value = 6302.9997
$ws.Range("M122").Value = -3852.9997

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 5736.125
$ws.Range("J80").Value = 9998.5
$ws.Range("L80").Value = 9998.5
$ws.Range("N80").Value = -11994.5
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 5736.125
$ws.Range("J83").Value = 9998.5
$ws.Range("L83").Value = 49992.5
$ws.Range("N83").Value = -59976.5
# Row 93: One Ring Circus
$ws.Range("H93").Value = 59999
$ws.Range("J93").Value = 59999
$ws.Range("L93").Value = 59999
$ws.Range("N93").Value = -63743
# Row 109: You're My Wonderhall
$ws.Range("H109").Value = 44751.4
$ws.Range("J109").Value = 44751.4
$ws.Range("L109").Value = 44751.4
$ws.Range("N109").Value = -46831.4
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 4052.4119
$ws.Range("I126").Value = 3662.818
$ws.Range("K126").Value = 10988.454
$ws.Range("M126").Value = -8518.454000000002
# Row 132: On Board for Lar
$ws.Range("H132").Value = 23258168
$ws.Range("I132").Value = 28573920
$ws.Range("K132").Value = 85721760
$ws.Range("M132").Value = -85719230

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 3138.3333
$ws.Range("I40").Value = 2733.5454
$ws.Range("J40").Value = 4251.5
$ws.Range("K40").Value = 2733.5454
$ws.Range("L40").Value = 4251.5
$ws.Range("M40").Value = -2597.5454
$ws.Range("N40").Value = -4523.5
# Row 122: Hell on Leather
$ws.Range("H122").Value = 5917
$ws.Range("I122").Value = 5558.294
$ws.Range("K122").Value = 16674.882
$ws.Range("M122").Value = -14224.882
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 112831.586
$ws.Range("I136").Value = 22000
$ws.Range("K136").Value = 66000
$ws.Range("M136").Value = -63450

$ws = $wb.Worksheets.Item("WVR")
# Row 70: An Account of My Boots
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
# Row 73: Soot in My Hair and Scars on My Feet (L)
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
